# h2 f3 error correction
# Update the enquiry date stamp and a handful of measured values in the
# technical-specification tables.

$d = $word.ActiveDocument

# --- Table 1 (header block): Enquiry Date ---------------------------------
# Row 2 ("Enquiry" | "a" | "Date" | <value>)
$d.Tables.Item(1).Cell(2, 4).Range.Text = "12/20/2019, 04:40 PM"

# --- Table 4 (B. COOLING WATER CIRCUIT) ------------------------------------
# Row 2: Cooling water flow (m3/hr)              114   -> 125
$d.Tables.Item(4).Cell(2, 4).Range.Text = "125"
# Row 4: Cooling water outlet temperature (C)     37.1 -> 36.6
$d.Tables.Item(4).Cell(4, 4).Range.Text = "36.6"
# Row 7: Cooling water circuit pressure loss (mLC) 2.2 -> 2.6
$d.Tables.Item(4).Cell(7, 4).Range.Text = "2.6"

# --- Table 5 (C. Steam Circuit) --------------------------------------------
# Row 3: Steam Consumption(+/-3%) (kg/hr)        400.5 -> 398.1
$d.Tables.Item(5).Cell(3, 4).Range.Text = "398.1"

# --- Table 6 (D. Electrical Data) ------------------------------------------
# Row 6: Vacuum pump rating (kW (A))        0.8( 1.8 ) -> 0.75( 1.8 )
$d.Tables.Item(6).Cell(6, 4).Range.Text = "0.75( 1.8 )"
